# Simulate the user dragging "Spinner 2" (the CHUNKSIZE spin button, linked to $B$6)
# from 4 down to 2. Excel form controls write their current position straight into
# the linked cell, which then ripples through the dependent formulas
# (C6/D6/C8:C11/D8:D11/E8:E11) via the normal recalculation engine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the spin button's own control properties (min/max/position) in sync too.
# "Drehfeld 2" / "Spinner 2" is the 2nd shape on the sheet and is the one whose
# fmlaLink points at $B$6 (the CHUNKSIZE spinner).
try {
    $spinner = $ws.Shapes.Item(2)
    $ctrl = $spinner.ControlFormat
    $ctrl.Value = 2
} catch {
    # Older/partial host implementations may not expose ControlFormat fully;
    # the linked-cell write below is what actually drives the worksheet.
}

# This is the authoritative change: the spinner's linked cell.
$ws.Range("B6").Value = 2

# The user had also moved the selection to B4 (row with the SIZE spinner) afterwards.
$ws.Range("B4").Select()

$wb.Application.CalculateFull()
